$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 70.818184
$ws.Range("I2").Value = 48
$ws.Range("K2").Value = 48
$ws.Range("M2").Value = 65

$ws.Range("H33").Value = 654.375
$ws.Range("I33").Value = 686.35
$ws.Range("K33").Value = 686.35
$ws.Range("M33").Value = -457.35

$ws.Range("H112").Value = 6120
$ws.Range("J112").Value = 6260.5713
$ws.Range("L112").Value = 18781.7139
$ws.Range("N112").Value = -20997.7139

$ws.Range("H138").Value = 4665.1665
$ws.Range("I138").Value = 1962.8823
$ws.Range("J138").Value = 6147.0645
$ws.Range("K138").Value = 5888.6469
$ws.Range("L138").Value = 18441.1935
$ws.Range("M138").Value = -748.6468999999997
$ws.Range("N138").Value = -28721.1935

$ws.Range("H141").Value = 1901.25
$ws.Range("I141").Value = 1833.3334
$ws.Range("K141").Value = 5500.0002
$ws.Range("M141").Value = -320.0002000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2099.0952
$ws.Range("I32").Value = 1840.0375
$ws.Range("J32").Value = 7280.25
$ws.Range("K32").Value = 1840.0375
$ws.Range("L32").Value = 7280.25
$ws.Range("M32").Value = -1553.0375
$ws.Range("N32").Value = -7854.25

$ws.Range("H61").Value = 5307.2827
$ws.Range("I61").Value = 2544.9312
$ws.Range("J61").Value = 10019.529
$ws.Range("K61").Value = 2544.9312
$ws.Range("L61").Value = 10019.529
$ws.Range("M61").Value = -2332.9312
$ws.Range("N61").Value = -10443.529

$ws.Range("H74").Value = 18039.25
$ws.Range("I74").Value = 23309.482
$ws.Range("J74").Value = 4145
$ws.Range("K74").Value = 23309.482
$ws.Range("L74").Value = 4145
$ws.Range("M74").Value = -22435.482
$ws.Range("N74").Value = -5893

$ws.Range("H77").Value = 18039.25
$ws.Range("I77").Value = 23309.482
$ws.Range("J77").Value = 4145
$ws.Range("K77").Value = 116547.41
$ws.Range("L77").Value = 20725
$ws.Range("M77").Value = -112179.41
$ws.Range("N77").Value = -29461

$ws.Range("H102").Value = 2419.8125
$ws.Range("I102").Value = 2407.5715
$ws.Range("K102").Value = 2407.5715
$ws.Range("M102").Value = -785.5715

$ws.Range("H122").Value = 11374.954
$ws.Range("I122").Value = 14458.6
$ws.Range("K122").Value = 43375.8
$ws.Range("M122").Value = -40925.8

$ws.Range("H132").Value = 6725.575
$ws.Range("I132").Value = 5314.2
$ws.Range("J132").Value = 9077.866
$ws.Range("K132").Value = 15942.6
$ws.Range("L132").Value = 27233.598
$ws.Range("M132").Value = -13412.6
$ws.Range("N132").Value = -32293.598

$ws.Range("H136").Value = 5307.2827
$ws.Range("I136").Value = 2544.9312
$ws.Range("J136").Value = 10019.529
$ws.Range("K136").Value = 7634.7936
$ws.Range("L136").Value = 30058.587
$ws.Range("M136").Value = -5084.7936
$ws.Range("N136").Value = -35158.587

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 31250324
$ws.Range("J80").Value = 295.55554
$ws.Range("L80").Value = 295.55554
$ws.Range("N80").Value = -2291.55554

$ws.Range("H83").Value = 31250324
$ws.Range("J83").Value = 295.55554
$ws.Range("L83").Value = 1477.7777
$ws.Range("N83").Value = -11461.7777

$ws.Range("H94").Value = 1736.1613
$ws.Range("I94").Value = 570.4
$ws.Range("K94").Value = 570.4
$ws.Range("M94").Value = -119.4

$ws.Range("H134").Value = 4308.677
$ws.Range("I134").Value = 1633.5641
$ws.Range("K134").Value = 4900.692300000001
$ws.Range("M134").Value = -2365.692300000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4082.2646
$ws.Range("I16").Value = 2448.0625
$ws.Range("J16").Value = 5534.8887
$ws.Range("K16").Value = 2448.0625
$ws.Range("L16").Value = 5534.8887
$ws.Range("M16").Value = -2161.0625
$ws.Range("N16").Value = -6108.8887

$ws.Range("H113").Value = 4082.2646
$ws.Range("I113").Value = 2448.0625
$ws.Range("J113").Value = 5534.8887
$ws.Range("K113").Value = 2448.0625
$ws.Range("L113").Value = 5534.8887
$ws.Range("M113").Value = -278.0625
$ws.Range("N113").Value = -9874.8887

$ws.Range("H122").Value = 1843
$ws.Range("I122").Value = 1790.6666
$ws.Range("K122").Value = 5371.9998
$ws.Range("M122").Value = -2921.9998

$ws.Range("H132").Value = 4828.4863
$ws.Range("I132").Value = 1574.8636
$ws.Range("J132").Value = 9600.467000000001
$ws.Range("K132").Value = 4724.5908
$ws.Range("L132").Value = 28801.401
$ws.Range("M132").Value = -2194.5908
$ws.Range("N132").Value = -33861.401

$ws.Range("H134").Value = 7393.738
$ws.Range("I134").Value = 6754.2188
$ws.Range("K134").Value = 20262.6564
$ws.Range("M134").Value = -17727.6564

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3080535
$ws.Range("I5").Value = 8001710.5
$ws.Range("J5").Value = 4800.375
$ws.Range("K5").Value = 24005131.5
$ws.Range("L5").Value = 14401.125
$ws.Range("M5").Value = -24005019.5
$ws.Range("N5").Value = -14625.125

$ws.Range("H12").Value = 2778582.5
$ws.Range("I12").Value = 811.1667
$ws.Range("J12").Value = 4167468
$ws.Range("K12").Value = 2433.5001
$ws.Range("L12").Value = 12502404
$ws.Range("M12").Value = -2260.5001
$ws.Range("N12").Value = -12502750

$ws.Range("H34").Value = 4575.9546
$ws.Range("J34").Value = 6223.5
$ws.Range("L34").Value = 18670.5
$ws.Range("N34").Value = -18838.5

$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws.Range("H135").Value = 3080535
$ws.Range("I135").Value = 8001710.5
$ws.Range("J135").Value = 4800.375
$ws.Range("K135").Value = 72015394.5
$ws.Range("L135").Value = 43203.375
$ws.Range("M135").Value = -72012859.5
$ws.Range("N135").Value = -48273.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3478.3845
$ws.Range("I80").Value = 3213
$ws.Range("K80").Value = 3213
$ws.Range("M80").Value = -2215

$ws.Range("H83").Value = 3478.3845
$ws.Range("I83").Value = 3213
$ws.Range("K83").Value = 16065
$ws.Range("M83").Value = -11073

$ws.Range("H132").Value = 8513.038
$ws.Range("I132").Value = 3119.3845
$ws.Range("K132").Value = 9358.1535
$ws.Range("M132").Value = -6828.1535

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7559.4707
$ws.Range("I7").Value = 6083.6665
$ws.Range("K7").Value = 6083.6665
$ws.Range("M7").Value = -5971.6665

$ws.Range("H16").Value = 1459.8966
$ws.Range("I16").Value = 1459.8966
$ws.Range("K16").Value = 1459.8966
$ws.Range("M16").Value = -1289.8966

$ws.Range("H46").Value = 7409722.5
$ws.Range("J46").Value = 7938953.5
$ws.Range("L46").Value = 7938953.5
$ws.Range("N46").Value = -7939329.5

$ws.Range("H55").Value = 111112050
$ws.Range("I55").Value = 1000000000
$ws.Range("J55").Value = 1050
$ws.Range("K55").Value = 1000000000
$ws.Range("L55").Value = 1050
$ws.Range("M55").Value = -999999827
$ws.Range("N55").Value = -1396

$ws.Range("H61").Value = 4888.125
$ws.Range("I61").Value = 2243.1
$ws.Range("J61").Value = 6777.4287
$ws.Range("K61").Value = 2243.1
$ws.Range("L61").Value = 6777.4287
$ws.Range("M61").Value = -2041.1
$ws.Range("N61").Value = -7181.4287

$ws.Range("H93").Value = 8583.416999999999
$ws.Range("I93").Value = 7875.375
$ws.Range("K93").Value = 7875.375
$ws.Range("M93").Value = -6627.375

$ws.Range("H113").Value = 4888.125
$ws.Range("I113").Value = 2243.1
$ws.Range("J113").Value = 6777.4287
$ws.Range("K113").Value = 2243.1
$ws.Range("L113").Value = 6777.4287
$ws.Range("M113").Value = -73.09999999999991
$ws.Range("N113").Value = -11117.4287

$ws.Range("H122").Value = 3724.2239
$ws.Range("I122").Value = 3239.7358
$ws.Range("K122").Value = 9719.207399999999
$ws.Range("M122").Value = -7269.207399999999

$ws.Range("H126").Value = 7559.4707
$ws.Range("I126").Value = 6083.6665
$ws.Range("K126").Value = 18250.9995
$ws.Range("M126").Value = -15780.9995

$ws.Range("H132").Value = 13166863
$ws.Range("I132").Value = 27783934
$ws.Range("J132").Value = 11499.95
$ws.Range("K132").Value = 83351802
$ws.Range("L132").Value = 34499.85000000001
$ws.Range("M132").Value = -83349272
$ws.Range("N132").Value = -39559.85000000001

$ws.Range("H136").Value = 11093.92
$ws.Range("I136").Value = 3119.6667
$ws.Range("J136").Value = 12181.318
$ws.Range("K136").Value = 9359.000100000001
$ws.Range("L136").Value = 36543.954
$ws.Range("M136").Value = -6809.000100000001
$ws.Range("N136").Value = -41643.954

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8341380.5
$ws.Range("I132").Value = 11908567
$ws.Range("J132").Value = 17945.334
$ws.Range("K132").Value = 35725701
$ws.Range("L132").Value = 53836.00199999999
$ws.Range("M132").Value = -35723171
$ws.Range("N132").Value = -58896.00199999999

$ws.Range("H136").Value = 29445754
$ws.Range("I136").Value = 111112380
$ws.Range("J136").Value = 45765.88
$ws.Range("K136").Value = 333337140
$ws.Range("L136").Value = 137297.64
$ws.Range("M136").Value = -333334590
$ws.Range("N136").Value = -142397.64
